$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddr, $val) {
    $c = $ws.Range($rangeAddr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue 'D2' '26.707.61'
Set-TextValue 'E2' '  -0.20%  '
Set-TextValue 'D3' '1.633.08'
Set-TextValue 'E3' '  -0.94%  '
Set-TextValue 'E4' '  +0.03%  '
Set-TextValue 'D5' '217.95'
Set-TextValue 'E5' '  +0.69%  '
Set-TextValue 'D6' '0.498'
Set-TextValue 'E6' '  -1.52%  '
Set-TextValue 'E7' '  +0.07%  '
Set-TextValue 'E8' '  -1.33%  '
Set-TextValue 'E9' '  -1.23%  '
Set-TextValue 'E10' '  -1.41%  '
Set-TextValue 'D11' '0.0843'
Set-TextValue 'E11' '  -0.03%  '
Set-TextValue 'D12' '1.860.85'
Set-TextValue 'D13' '1.626.55'
Set-TextValue 'E13' '  -2.14%  '
Set-TextValue 'E14' '  -2.17%  '
Set-TextValue 'E15' '  -2.06%  '
Set-TextValue 'D16' '63.94'
Set-TextValue 'E16' '  -2.05%  '
Set-TextValue 'D17' '26.673.09'
Set-TextValue 'E17' '  -0.30%  '
Set-TextValue 'E18' '  -3.36%  '
Set-TextValue 'B19' 'BitcoinCash'
Set-TextValue 'C19' 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue 'D19' '211.08'
Set-TextValue 'E19' '  -3.35%  '
Set-TextValue 'B20' 'Dai'
Set-TextValue 'C20' 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue 'D20' '1.01'
Set-TextValue 'E20' '  +0.05%  '
Set-TextValue 'D21' '4.30'
Set-TextValue 'E21' '  -1.68%  '
Set-TextValue 'B22' 'Chainlink'
Set-TextValue 'C22' 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue 'D22' '6.15'
Set-TextValue 'E22' '  -2.14%  '
Set-TextValue 'B23' 'Toncoin'
Set-TextValue 'C23' 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue 'D23' '2.32'
Set-TextValue 'E23' '  -8.85%  '
Set-TextValue 'E24' '  -3.17%  '
Set-TextValue 'D25' '146.59'
Set-TextValue 'E25' '  +0.05%  '
Set-TextValue 'E26' '  +0.30%  '
Set-TextValue 'E27' '  -2.50%  '
Set-TextValue 'D28' '7.01'
Set-TextValue 'E28' '  -2.40%  '
Set-TextValue 'D29' '15.51'
Set-TextValue 'E29' '  -1.78%  '
Set-TextValue 'D30' '0.0500'
Set-TextValue 'E30' '  -3.67%  '
Set-TextValue 'E31' '  +0.51%  '
Set-TextValue 'D32' '3.34'
Set-TextValue 'E32' '  -0.51%  '
Set-TextValue 'E33' '  -2.61%  '
Set-TextValue 'D34' '1.260.47'
Set-TextValue 'E34' '  -1.62%  '
Set-TextValue 'E35' '  +0.11%  '
Set-TextValue 'D36' '1.51'
Set-TextValue 'E36' '  -2.73%  '
Set-TextValue 'E37' '  -3.50%  '
Set-TextValue 'D38' '0.520'
Set-TextValue 'E38' '  -4.02%  '
Set-TextValue 'E39' '  +0.05%  '
Set-TextValue 'E40' '  -4.15%  '
Set-TextValue 'E41' '  -2.10%  '
Set-TextValue 'E42' '  -3.98%  '
Set-TextValue 'D43' '1.771.94'
Set-TextValue 'E43' '  -0.90%  '
Set-TextValue 'E44' '  -3.87%  '
Set-TextValue 'D45' '90.80'
Set-TextValue 'E45' '  -1.37%  '
Set-TextValue 'D46' '59.52'
Set-TextValue 'E46' '  -0.65%  '
Set-TextValue 'E47' '  -2.56%  '
Set-TextValue 'E48' '  +0.20%  '
Set-TextValue 'E49' '  +0.04%  '
Set-TextValue 'E50' '  -0.49%  '
Set-TextValue 'D51' '0.0954'
Set-TextValue 'E51' '  -2.45%  '
